$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Overview sheet: status text updated ("Ready for handoff" -> handback)
#    E2/F2 share the same string used elsewhere (Status columns on the
#    zh-cn / de-de sheets), so updating their text updates all of them.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

$wb.Worksheets.Item("zh-cn").Range("C2").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("de-de").Range("C2").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de status columns on the Overview sheet.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# Helper info reused for both language sheets
# ---------------------------------------------------------------------
$mdFileName = "69f57c00-3142-43d1-8e83-b210ad9ab90c.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b498c8a826b3da58c039a7f17d0c3a35c8d7cb9a/e2e/69f57c00-3142-43d1-8e83-b210ad9ab90c.md"

# ---------------------------------------------------------------------
# 2) zh-cn sheet: fill in the "Latest Target File" / "Latest Handback
#    File" / "Latest Handback DateTime" columns for the handback.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFileName)
$wsZhCn.Range("I2").Font.Name = "Calibri"
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = 15570276

$wsZhCn.Range("J2").Value = "69f57c00-3142-43d1-8e83-b210ad9ab90c.8407c903d4f40f0933fdbca2d05fba4e50c0eb26.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-02 13:08:34"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# 3) de-de sheet: same treatment as zh-cn, but a later handback
#    timestamp, so it picks up a brand-new datetime string.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFileName)
$wsDeDe.Range("I2").Font.Name = "Calibri"
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = 15570276

$wsDeDe.Range("J2").Value = "69f57c00-3142-43d1-8e83-b210ad9ab90c.8407c903d4f40f0933fdbca2d05fba4e50c0eb26.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-02 13:08:42"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17
